# Adicionados balancos concatenados em uma unica planilha.
# Adds a new column AH containing the 30/06/2024 balance sheet figures,
# mirroring the existing column layout (one column per reporting date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = 34  # column AH

# 1) Header cell AH1: copy formatting (style) from AG1, then set the new date text.
$ws.Range("AG1").Copy($ws.Cells.Item(1, $col))
$ws.Cells.Item(1, $col).Value = "30/06/2024"

# 2) Numeric data rows: row number -> value for 30/06/2024
$values = @{
    2 = 12851.859
    3 = 2188.061
    4 = 47.879
    5 = 1119.174
    6 = 594.175
    7 = 0
    8 = 0
    9 = 101.518
    10 = 0
    11 = 325.315
    12 = 803.569
    13 = 0
    14 = 0
    15 = 269.354
    16 = 19.262
    17 = 0
    18 = 0
    19 = 3.569
    20 = 0
    21 = 54.375
    22 = 53.352
    23 = 7942.297
    24 = 1864.58
    25 = 0
    26 = 12851.859
    27 = 1177.842
    28 = 23.808
    29 = 222.036
    30 = 74.896
    31 = 625.915
    32 = 27.794
    33 = 0
    34 = 203.393
    35 = 0
    36 = 0
    37 = 7021.378
    38 = 5905.636
    39 = 0
    40 = 575.013
    41 = 540.729
    42 = 0
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 4652.639
    48 = 3968.679
    49 = 62.54
    50 = 0
    51 = 804.378
    52 = 205.576
    53 = -388.534
    54 = 0
    55 = 0
    56 = 0
    59 = 658.902
    60 = -515.924
    61 = 142.978
    62 = 0
    63 = -38.416
    64 = 0
    65 = -3.736
    66 = 0
    67 = 3.462
    68 = -135.587
    69 = 25.238
    70 = -160.825
    74 = -31.299
    75 = -16.945
    76 = 2.109
    79 = 0
    80 = -46.135
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, $col).Value = $values[$row]
}

# 3) Section-header / blank rows: these rows only carry an (empty) inline
#    string in every other column, so replicate that by copying the
#    already-blank AG cell of the same row into AH.
$blankRows = @(57, 58, 71, 72, 73, 77, 78)
foreach ($row in $blankRows) {
    $ws.Cells.Item($row, 33).Copy($ws.Cells.Item($row, $col))
}
